$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (7-11) appended to the IDAN dataset.
# Columns: A=course, C=date, D=time to tablet (minutes), E=class length (hours), H=break time (minutes)
# Formats for C and E are copied from an existing data row so no new
# number-format styles get introduced.

# Row 7: algorithems
$ws.Range("A7").Value = "algorithems"
$ws.Range("C2").Copy()
$ws.Range("C7").PasteSpecial(-4122)
$ws.Range("C7").Value = 45798
$ws.Range("D7").Value = 110
$ws.Range("E2").Copy()
$ws.Range("E7").PasteSpecial(-4122)
$ws.Range("E7").Value = 2.5 / 24
$ws.Range("H7").Value = 30

# Row 8: architecture
$ws.Range("A8").Value = "architecture"
$ws.Range("C2").Copy()
$ws.Range("C8").PasteSpecial(-4122)
$ws.Range("C8").Value = 45798
$ws.Range("D8").Value = 8
$ws.Range("E2").Copy()
$ws.Range("E8").PasteSpecial(-4122)
$ws.Range("E8").Value = 1.5 / 24
$ws.Range("H8").Value = 0

# Row 9: algorithems
$ws.Range("A9").Value = "algorithems"
$ws.Range("C2").Copy()
$ws.Range("C9").PasteSpecial(-4122)
$ws.Range("C9").Value = 45799
$ws.Range("D9").Value = 0
$ws.Range("E2").Copy()
$ws.Range("E9").PasteSpecial(-4122)
$ws.Range("E9").Value = 1.5 / 24

# Row 10: algorithems
$ws.Range("A10").Value = "algorithems"
$ws.Range("C2").Copy()
$ws.Range("C10").PasteSpecial(-4122)
$ws.Range("C10").Value = 45799
$ws.Range("D10").Value = 90
$ws.Range("E2").Copy()
$ws.Range("E10").PasteSpecial(-4122)
$ws.Range("E10").Value = 1.5 / 24

# Row 11: kaplt
$ws.Range("A11").Value = "kaplt"
$ws.Range("C2").Copy()
$ws.Range("C11").PasteSpecial(-4122)
$ws.Range("C11").Value = 45799
$ws.Range("D11").Value = 120
$ws.Range("E2").Copy()
$ws.Range("E11").PasteSpecial(-4122)
$ws.Range("E11").Value = 2.0 / 24
$ws.Range("H11").Value = 15

$excel.CutCopyMode = 0

# Update the selection to reflect the cursor location after the edit.
$ws.Range("B26").Select()
